{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the CV content updates described by the commit:\n//  - New home address (Chula Vista, CA instead of Hong Kong)\n//  - New US mobile number\n//  - Remove the stray \"_GoBack\" bookmark\n//  - Updated job title / company name for the LF Logistics role\n//  - \"Li & Fung's\" -> \"LF's\" wording tweaks\n//  - IPG Mediabrands region updated from (Philippines) to (Asia Pacific)\n\nasync function replaceAll(searchRange, searchText, newText) {\n  const results = searchRange.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nasync function replaceWithinParagraphOf(anchorText, searchText, newText) {\n  const results = context.document.body.search(anchorText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Anchor text not found: \" + anchorText);\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  await replaceAll(para, searchText, newText);\n}\n\nconst body = context.document.body;\n\n// 1) Residence address (the whole old address is being replaced).\nawait replaceAll(\n  body,\n  \"Flat C 16th Floor Tower 10, Park Central, 9 Tong Tak Street, Tseung Kwan O\",\n  \"1082 Mt Dana Dr, Chula Vista, CA 91913\"\n);\n\n// 2) Mobile number.\nawait replaceAll(body, \"+852.9732.6715\", \"+1-619-800-0859\\u202C\");\n\n// 3) Remove the leftover \"_GoBack\" bookmark (South China Morning Post sentence).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 4) Job title for the LF Logistics role (whole line is replaced).\nawait replaceAll(\n  body,\n  \"Data Science Manager, Global Supply Chain Analytics\",\n  \"Senior Manager \u2013 Data Science, Global Supply Chain Analytics\"\n);\n\n// 5) Company name/location line for the LF Logistics role. Keep the leading\n//    \"LF Logistics \" run untouched and only rewrite the trailing portion.\nawait replaceWithinParagraphOf(\n  \"LF Logistics\",\n  \"\u2013 A Li & Fung Company (Hong Kong)\",\n  \"(Global)\"\n);\n\n// 6) \"Li & Fung's\" -> \"LF's\" (two bullet points) \u2014 narrow substring replace\n//    so surrounding runs/sentences stay untouched.\nawait replaceAll(body, \"Li & Fung\\u2019s\", \"LF\\u2019s\");\n\n// 7) IPG Mediabrands region \u2014 only rewrite the \"(Philippines)\" run, scoped\n//    to that paragraph so the earlier \"IPG \"/\"Mediabrands\" runs (and their\n//    spell-check markers) stay untouched.\nawait replaceWithinParagraphOf(\"IPG Mediabrands\", \" (Philippines)\", \" (Asia Pacific)\");\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the CV content updates described by the commit:\n#  - New home address (Chula Vista, CA instead of Hong Kong)\n#  - New US mobile number\n#  - Remove the stray \"_GoBack\" bookmark\n#  - Updated job title / company name for the LF Logistics role\n#  - \"Li & Fung's\" -> \"LF's\" wording tweaks\n#  - IPG Mediabrands region updated from (Philippines) to (Asia Pacific)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Find text not found: $findText\"\n    }\n}\n\n# 1) Residence address.\nReplace-Text \"Flat C 16th Floor Tower 10, Park Central, 9 Tong Tak Street, Tseung Kwan O\" \"1082 Mt Dana Dr, Chula Vista, CA 91913\"\n\n# 2) Mobile number (includes a trailing pop-directional-formatting mark, as in the source).\nReplace-Text \"+852.9732.6715\" \"+1-619-800-0859\u202c\"\n\n# 3) Remove the leftover \"_GoBack\" bookmark (South China Morning Post sentence).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 4) Job title for the LF Logistics role.\nReplace-Text \"Data Science Manager, Global Supply Chain Analytics\" \"Senior Manager \u2013 Data Science, Global Supply Chain Analytics\"\n\n# 5) Company name/location line for the LF Logistics role.\nReplace-Text \"LF Logistics \u2013 A Li & Fung Company (Hong Kong)\" \"LF Logistics (Global)\"\n\n# 6) \"Li & Fung's\" -> \"LF's\" (two bullet points).\nReplace-Text \"Responsible for optimizing Li & Fung\u2019s global supply chain operations\" \"Responsible for optimizing LF\u2019s global supply chain operations\"\nReplace-Text \"Developed the data science platform for Li & Fung\u2019s costing center of excellence\" \"Developed the data science platform for LF\u2019s costing center of excellence\"\n\n# 7) IPG Mediabrands region.\nReplace-Text \"IPG Mediabrands (Philippines)\" \"IPG Mediabrands (Asia Pacific)\"\n"}
